# feat(products/import): improve field mapping, UI labels, and documentation
# for product schema.
#
# Replaces the English/technical machine-readable column headers on the
# sample import sheet with the localized Russian labels that the product
# catalog UI now uses, and normalizes the barcode column so it stores plain
# numbers instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old header -> new localized header, in column order (A..AB).
$headers = @(
    "SKU",                        # sku
    "Артикул продавца",           # seller_sku
    "Артикул WB",                 # wb_sku
    "NM ID",                      # nm_id
    "Название",                   # title
    "Бренд",                      # brand
    "Категория",                  # category
    "Текущая цена",               # price_src
    "Скидка, %",                  # seller_discount_pct
    "Итоговая цена (расчёт)",     # price
    "Цена со скидкой (расчёт)",   # price_final
    "Остаток общий (расчёт)",     # stock
    "Остатки WB",                 # stock_wb
    "Остатки продавца",           # stock_seller
    "Штрихкод",                   # barcode
    "Активен",                    # is_active
    "Себик",                      # product_cost
    "Транспортировка",            # shipping_cost
    "Логистика возврата",         # logistics_back_cost
    "Коэфф. склада",              # warehouse_coeff
    "Оборачиваемость, дни",       # turnover_days
    "Вес с упаковкой (кг)",       # weight_kg
    "Длина упаковки, см",         # package_l_cm
    "Ширина упаковки, см",        # package_w_cm
    "Высота упаковки, см",        # package_h_cm
    "Литраж",                     # volume_l
    "Комменты",                   # comments
    "Доп данные JSON"             # custom_data
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The barcode column (O) was imported as text; store it as a genuine number
# for every sample data row so downstream numeric handling/import works.
$lastRow = 4
$barcodeCol = 15
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $barcodeCol)
    $cell.Value = [double]$cell.Value2
}
